# Applies the "Automatic update of files" edit to the artfynd workbook.
# The underlying observation records got re-matched to different rows;
# this script rewrites the affected rows' data cells in place so that
# each row ends up holding the record the diff says it should hold,
# while row-invariant columns (site name, accuracy, county, dates, ...)
# stay untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Record($Row, $A, $B, $D, $E, $F, $G, $H, $Q, $R) {
    $ws.Cells.Item($Row, 1).Value2  = $A   # A - Id
    $ws.Cells.Item($Row, 2).Value2  = $B   # B - Taxonsorteringsordning
    $ws.Cells.Item($Row, 4).Value2  = $D   # D - Rödlistade
    $ws.Cells.Item($Row, 5).Value2  = $E   # E - TaxonId
    $ws.Cells.Item($Row, 6).Value2  = $F   # F - Artnamn
    $ws.Cells.Item($Row, 7).Value2  = $G   # G - Vetenskapligt namn
    $ws.Cells.Item($Row, 8).Value2  = $H   # H - Auktor
    $ws.Cells.Item($Row, 17).Value2 = $Q   # Q - Ost
    $ws.Cells.Item($Row, 18).Value2 = $R   # R - Nord
}

# --- Rows 6 and 8 swap their species/coordinate data -----------------
Set-Record 6 131066787 80384 "LC" 6463 "Bårdlav" "Nephroma parile" "(Ach.) Ach." 425069 6712290
Set-Record 8 131066788 83216 "NT" 308 "Brunpudrad nållav" "Chaenotheca gracillima" "(Vain.) Tibell" 425211 6712276

# --- Rows 15, 16, 17 rotate their species/coordinate data -------------
Set-Record 15 131066768 91809 "NT" 1202 "Ullticka" "Phellinidium ferrugineofuscum" "(P.Karst.) Fiasson & Niemelä" 425256 6712203
Set-Record 16 131066761 91772 "LC" 5447 "Vedticka" "Fuscoporia viticola" "(Schwein.) Murrill" 425072 6712273
Set-Record 17 131066782 91823 "NT" 1204 "Gränsticka" "Phellopilus nigrolimitatus" "(Romell) Niemelä, T.Wagner & M.Fisch." 425059 6712253

# --- Rows 19 and 20 swap their Id / coordinates (same species) --------
$ws.Cells.Item(19, 1).Value2  = 131066775
$ws.Cells.Item(19, 17).Value2 = 425244
$ws.Cells.Item(19, 18).Value2 = 6712292

$ws.Cells.Item(20, 1).Value2  = 131066773
$ws.Cells.Item(20, 17).Value2 = 425271
$ws.Cells.Item(20, 18).Value2 = 6712264

# --- Rows 21 and 23 fully swap their records, including the optional --
# --- Ålder-Stadium/Aktivitet (K,L,M,N) and Publik kommentar (AC) cells -
Set-Record 21 131066778 81229 "NT" 1049 "Kortskaftad ärgspik" "Microcalicium ahlneri" "Tibell" 425336 6712202
$ws.Cells.Item(21, 11).ClearContents()  # K21
$ws.Cells.Item(21, 12).ClearContents()  # L21
$ws.Cells.Item(21, 13).ClearContents()  # M21 (äldre spår)
$ws.Cells.Item(21, 14).ClearContents()  # N21
$ws.Cells.Item(21, 29).ClearContents()  # AC21 (Ringhack på gran)

Set-Record 23 131066774 57884 "NT" 100109 "Tretåig hackspett" "Picoides tridactylus" "(Linnaeus, 1758)" 425250 6712265
# K23, L23, N23 need to exist as present-but-empty cells (matching the
# "äldre spår" rows elsewhere in the sheet). Plain Value2 assignment of ""
# does not materialize an empty cell, so copy the already-empty I23 cell
# into them instead.
$ws.Cells.Item(23, 9).Copy($ws.Cells.Item(23, 11))  # K23
$ws.Cells.Item(23, 9).Copy($ws.Cells.Item(23, 12))  # L23
$ws.Cells.Item(23, 13).Value2 = "äldre spår"        # M23
$ws.Cells.Item(23, 9).Copy($ws.Cells.Item(23, 14))  # N23
$ws.Cells.Item(23, 29).Value2 = "Ringhack på gran"  # AC23

# --- Rows 32, 33, 34 rotate their species/coordinate data --------------
Set-Record 32 131066784 83090 "NT" 1312 "Gammelgransskål" "Pseudographis pinicola" "(Nyl.) Rehm" 425170 6712283
Set-Record 33 131066767 91809 "NT" 1202 "Ullticka" "Phellinidium ferrugineofuscum" "(P.Karst.) Fiasson & Niemelä" 425259 6712201
Set-Record 34 131066790 83216 "NT" 308 "Brunpudrad nållav" "Chaenotheca gracillima" "(Vain.) Tibell" 425164 6712278

$wb.Save()
